# Add I0 and IF columns (I and J) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Apply header style consistent with other header cells (B1:H1)
$ws.Range("I1:J1").Style = $ws.Range("H1").Style

# Data values for I and J columns, rows 2..21
$values = @(
    @(8, 8),
    @(4, 4),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(6, 7),
    @(7, 7),
    @(4, 5),
    @(7, 7),
    @(9, 9),
    @(6, 6),
    @(4, 5),
    @(5, 6),
    @(6, 7),
    @(8, 8),
    @(7, 7),
    @(7, 7),
    @(3, 4),
    @(7, 7),
    @(8, 8)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
